$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("L69").ClearContents()
$ws.Range("N69").Value2 = 0
$ws.Range("H72").Value2 = 0
$ws.Range("J72").Value2 = 0
$ws.Range("L72").ClearContents()
$ws.Range("N72").Value2 = 0
$ws.Range("H101").Value2 = 547.7778
$ws.Range("I101").Value2 = 520.2857
$ws.Range("K101").Value2 = 1560.8571
$ws.Range("M101").Value2 = 61.14289999999983
$ws.Range("H107").Value2 = 730.381
$ws.Range("I107").Value2 = 804.6111
$ws.Range("J107").Value2 = 285
$ws.Range("K107").Value2 = 804.6111
$ws.Range("L107").Value2 = 285
$ws.Range("M107").Value2 = 1115.3889
$ws.Range("N107").Value2 = -4125
$ws.Range("H125").Value2 = 26109.076
$ws.Range("I125").Value2 = 65195.8
$ws.Range("J125").Value2 = 1679.875
$ws.Range("K125").Value2 = 586762.2000000001
$ws.Range("L125").Value2 = 15118.875
$ws.Range("M125").Value2 = -584302.2000000001
$ws.Range("N125").Value2 = -20038.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 3465.7
$ws.Range("I2").Value2 = 3699.5
$ws.Range("J2").Value2 = 3231.9
$ws.Range("K2").Value2 = 3699.5
$ws.Range("L2").Value2 = 3231.9
$ws.Range("M2").Value2 = -3586.5
$ws.Range("N2").Value2 = -3457.9
$ws.Range("H32").Value2 = 44268.54
$ws.Range("I32").Value2 = 25568.238
$ws.Range("J32").Value2 = 142445.12
$ws.Range("K32").Value2 = 25568.238
$ws.Range("L32").Value2 = 142445.12
$ws.Range("M32").Value2 = -25281.238
$ws.Range("N32").Value2 = -143019.12
$ws.Range("H74").Value2 = 1414.228
$ws.Range("I74").Value2 = 1310.551
$ws.Range("J74").Value2 = 2049.25
$ws.Range("K74").Value2 = 1310.551
$ws.Range("L74").Value2 = 2049.25
$ws.Range("M74").Value2 = -436.5509999999999
$ws.Range("N74").Value2 = -3797.25
$ws.Range("H77").Value2 = 1414.228
$ws.Range("I77").Value2 = 1310.551
$ws.Range("J77").Value2 = 2049.25
$ws.Range("K77").Value2 = 6552.754999999999
$ws.Range("L77").Value2 = 10246.25
$ws.Range("M77").Value2 = -2184.754999999999
$ws.Range("N77").Value2 = -18982.25
$ws.Range("H101").Value2 = 28637.2
$ws.Range("J101").Value2 = 28637.2
$ws.Range("L101").Value2 = 28637.2
$ws.Range("N101").Value2 = -35127.2
$ws.Range("H116").Value2 = 3465.7
$ws.Range("I116").Value2 = 3699.5
$ws.Range("J116").Value2 = 3231.9
$ws.Range("K116").Value2 = 3699.5
$ws.Range("L116").Value2 = 3231.9
$ws.Range("M116").Value2 = -1405.5
$ws.Range("N116").Value2 = -7819.9
$ws.Range("H125").Value2 = 55428.4
$ws.Range("J125").Value2 = 55428.4
$ws.Range("L125").Value2 = 55428.4
$ws.Range("N125").Value2 = -65268.4
$ws.Range("H127").Value2 = 0
$ws.Range("J127").Value2 = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value2 = 0
$ws.Range("H128").Value2 = 77263.5
$ws.Range("J128").Value2 = 77263.5
$ws.Range("L128").Value2 = 77263.5
$ws.Range("N128").Value2 = -87223.5
$ws.Range("H129").Value2 = 100680
$ws.Range("J129").Value2 = 100680
$ws.Range("L129").Value2 = 100680
$ws.Range("N129").Value2 = -110680
$ws.Range("H132").Value2 = 2111.2812
$ws.Range("I132").Value2 = 1270.3214
$ws.Range("K132").Value2 = 3810.9642
$ws.Range("M132").Value2 = -1280.9642
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 3465.7
$ws.Range("I3").Value2 = 3699.5
$ws.Range("J3").Value2 = 3231.9
$ws.Range("K3").Value2 = 3699.5
$ws.Range("L3").Value2 = 3231.9
$ws.Range("M3").Value2 = -3585.5
$ws.Range("N3").Value2 = -3459.9
$ws.Range("H94").Value2 = 567.85187
$ws.Range("I94").Value2 = 567.85187
$ws.Range("K94").Value2 = 567.85187
$ws.Range("M94").Value2 = -116.85187
$ws.Range("H105").Value2 = 4826.85
$ws.Range("I105").Value2 = 5112.5557
$ws.Range("J105").Value2 = 2255.5
$ws.Range("K105").Value2 = 5112.5557
$ws.Range("L105").Value2 = 2255.5
$ws.Range("M105").Value2 = -3365.5557
$ws.Range("N105").Value2 = -5749.5
$ws.Range("H129").Value2 = 66678.5
$ws.Range("J129").Value2 = 66678.5
$ws.Range("L129").Value2 = 66678.5
$ws.Range("N129").Value2 = -76678.5
$ws.Range("H130").Value2 = 51194.75
$ws.Range("J130").Value2 = 51194.75
$ws.Range("L130").Value2 = 51194.75
$ws.Range("N130").Value2 = -61234.75
$ws.Range("H134").Value2 = 1746.55
$ws.Range("I134").Value2 = 1759.5264
$ws.Range("J134").Value2 = 1500
$ws.Range("K134").Value2 = 5278.5792
$ws.Range("L134").Value2 = 4500
$ws.Range("M134").Value2 = -2743.5792
$ws.Range("N134").Value2 = -9570
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1583.1111
$ws.Range("I16").Value2 = 1291.5
$ws.Range("K16").Value2 = 1291.5
$ws.Range("M16").Value2 = -1004.5
$ws.Range("H31").Value2 = 3215.1924
$ws.Range("I31").Value2 = 1776.7
$ws.Range("K31").Value2 = 1776.7
$ws.Range("M31").Value2 = -1481.7
$ws.Range("H34").Value2 = 3215.1924
$ws.Range("I34").Value2 = 1776.7
$ws.Range("K34").Value2 = 1776.7
$ws.Range("M34").Value2 = -1574.7
$ws.Range("H62").Value2 = 4666.5557
$ws.Range("J62").Value2 = 0
$ws.Range("L62").Value2 = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value2 = 4666.5557
$ws.Range("J65").Value2 = 0
$ws.Range("L65").Value2 = 0
$ws.Range("N65").ClearContents()
$ws.Range("H87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value2 = 0
$ws.Range("H90").Value2 = 0
$ws.Range("J90").Value2 = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value2 = 0
$ws.Range("H107").Value2 = 1471.8889
$ws.Range("J107").Value2 = 2187.1428
$ws.Range("L107").Value2 = 2187.1428
$ws.Range("N107").Value2 = -6027.1428
$ws.Range("H113").Value2 = 1583.1111
$ws.Range("I113").Value2 = 1291.5
$ws.Range("K113").Value2 = 1291.5
$ws.Range("M113").Value2 = 878.5
$ws.Range("H132").Value2 = 1637.2094
$ws.Range("I132").Value2 = 1590.2439
$ws.Range("K132").Value2 = 4770.7317
$ws.Range("M132").Value2 = -2240.7317
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value2 = 8987.888999999999
$ws.Range("I56").Value2 = 8987.888999999999
$ws.Range("K56").Value2 = 8987.888999999999
$ws.Range("M56").Value2 = -8457.888999999999
$ws.Range("H129").Value2 = 145836.58
$ws.Range("J129").Value2 = 3901.7778
$ws.Range("L129").Value2 = 11705.3334
$ws.Range("N129").Value2 = -21705.3334
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 34414.95
$ws.Range("I97").Value2 = 55494.082
$ws.Range("J97").Value2 = 2796.25
$ws.Range("K97").Value2 = 55494.082
$ws.Range("L97").Value2 = 2796.25
$ws.Range("M97").Value2 = -54998.082
$ws.Range("N97").Value2 = -3788.25
$ws.Range("H113").Value2 = 2399.3572
$ws.Range("I113").Value2 = 2628.1428
$ws.Range("J113").Value2 = 2170.5715
$ws.Range("K113").Value2 = 2628.1428
$ws.Range("L113").Value2 = 2170.5715
$ws.Range("M113").Value2 = -458.1428000000001
$ws.Range("N113").Value2 = -6510.5715
$ws.Range("I132").Value2 = 2447.4614
$ws.Range("J132").Value2 = 3588
$ws.Range("K132").Value2 = 7342.3842
$ws.Range("L132").Value2 = 10764
$ws.Range("M132").Value2 = -4812.3842
$ws.Range("N132").Value2 = -15824
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 3402.5
$ws.Range("I40").Value2 = 2771.3333
$ws.Range("J40").Value2 = 4664.8335
$ws.Range("K40").Value2 = 2771.3333
$ws.Range("L40").Value2 = 4664.8335
$ws.Range("M40").Value2 = -2635.3333
$ws.Range("N40").Value2 = -4936.8335
$ws.Range("H46").Value2 = 3138.6086
$ws.Range("I46").Value2 = 1807.3334
$ws.Range("J46").Value2 = 4590.909
$ws.Range("K46").Value2 = 1807.3334
$ws.Range("L46").Value2 = 4590.909
$ws.Range("M46").Value2 = -1619.3334
$ws.Range("N46").Value2 = -4966.909
$ws.Range("H61").Value2 = 51260.65
$ws.Range("I61").Value2 = 56723.555
$ws.Range("K61").Value2 = 56723.555
$ws.Range("M61").Value2 = -56521.555
$ws.Range("H113").Value2 = 51260.65
$ws.Range("I113").Value2 = 56723.555
$ws.Range("K113").Value2 = 56723.555
$ws.Range("M113").Value2 = -54553.555
$ws.Range("H122").Value2 = 10277.5625
$ws.Range("I122").Value2 = 13627.1
$ws.Range("J122").Value2 = 4695
$ws.Range("K122").Value2 = 40881.3
$ws.Range("L122").Value2 = 14085
$ws.Range("M122").Value2 = -38431.3
$ws.Range("N122").Value2 = -18985
$ws.Range("H132").Value2 = 2940.2156
$ws.Range("I132").Value2 = 2034.75
$ws.Range("K132").Value2 = 6104.25
$ws.Range("M132").Value2 = -3574.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value2 = 20004.166
$ws.Range("H73").Value2 = 20004.166
$ws.Range("H107").Value2 = 50001160
$ws.Range("I107").Value2 = 1201
$ws.Range("K107").Value2 = 3603
$ws.Range("M107").Value2 = -1683
$ws.Range("H113").Value2 = 1075.5294
$ws.Range("I113").Value2 = 509.44446
$ws.Range("K113").Value2 = 1528.33338
$ws.Range("M113").Value2 = 641.66662
$ws.Range("H126").Value2 = 5238.7915
$ws.Range("I126").Value2 = 2518
$ws.Range("J126").Value2 = 10680.375
$ws.Range("K126").Value2 = 7554
$ws.Range("L126").Value2 = 32041.125
$ws.Range("M126").Value2 = -5084
$ws.Range("N126").Value2 = -36981.125
$ws.Range("H132").Value2 = 3836.6978
$ws.Range("I132").Value2 = 3912.5715
$ws.Range("K132").Value2 = 11737.7145
$ws.Range("M132").Value2 = -9207.7145
